$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 8
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 8
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 8
$ws.Range("B10").Value = "<sentence>"
$ws.Range("C11").Value = 6
$ws.Range("C13").Value = 9
$ws.Range("C15").Value = 7
$ws.Range("C16").Value = 9
$ws.Range("C18").Value = 12
